# Release update: bump version to 0.1.1, refresh the publish date, add a
# "Jurisdiction" metadata row, and rename the second sheet's tab.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Insert a new row for "Jurisdiction" right above "Description" (row 11),
# copying formatting from the row above ("Contact") so the new row matches
# the existing body-row style.
$meta.Rows.Item(11).Insert()
$meta.Range("A10:B10").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)

$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# Bump the recorded Version.
$meta.Range("B3").Value = "0.1.1"

# Refresh the recorded Date to match the new release.
$meta.Range("B8").Value = "2024-11-11T17:53:38-06:00"

# Rename the "Include ..." sheet's tab to "Include #0".
$include = $wb.Worksheets.Item("Include from NMDP Language Co")
$include.Name = "Include #0"
